$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like pure numbers but must stay text (as in the source data)
# Pre-format them as Text so COM does not coerce the assigned string into a Double.
$textCells = @("D5","D6","D7","D9","D10","D11","D12","D13","D16","D17","D19","D21","D22","D23","D26","D27","D29","D30","D31","D32","D33","D35","D37","D38","D39","D42","D44","D45","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '46.032.05'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '2.360.28'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '302.23'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('D6').Value = '99.02'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '0.569'
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  -2.94%  '
$ws.Range('D10').Value = '34.55'
$ws.Range('E10').Value = '  -2.89%  '
$ws.Range('D11').Value = '0.0801'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').Value = '7.15'
$ws.Range('E12').Value = '  -2.87%  '
$ws.Range('D13').Value = '0.103'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '2.725.51'
$ws.Range('E14').Value = '  +2.33%  '
$ws.Range('D15').Value = '2.368.98'
$ws.Range('E15').Value = '  +2.38%  '
$ws.Range('D16').Value = '0.812'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').Value = '13.66'
$ws.Range('E17').Value = '  -2.27%  '
$ws.Range('D18').Value = '45.963.02'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('D19').Value = '12.85'
$ws.Range('E19').Value = '  -2.97%  '
$ws.Range('D20').Value = '0.0₃0972'
$ws.Range('E20').Value = '  +3.24%  '
$ws.Range('D21').Value = '6.04'
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range('D22').Value = '67.19'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').Value = '245.59'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('E24').Value = '  -2.78%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').Value = '1.93'
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('D27').Value = '39.97'
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('D29').Value = '9.83'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = '3.77'
$ws.Range('E30').Value = '  +19.94%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '20.99'
$ws.Range('E31').Value = '  +4.32%  '
$ws.Range('D32').Value = '2.77'
$ws.Range('E32').Value = '  +5.23%  '
$ws.Range('D33').Value = '5.55'
$ws.Range('E33').Value = '  -4.23%  '
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('D35').Value = '0.0777'
$ws.Range('E35').Value = '  -3.16%  '
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('D37').Value = '1.89'
$ws.Range('E37').Value = '  +4.64%  '
$ws.Range('D38').Value = '0.116'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('D39').Value = '15.17'
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('E41').Value = '  -1.91%  '
$ws.Range('D42').Value = '3.22'
$ws.Range('E42').Value = '  -7.14%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.882.75'
$ws.Range('E43').Value = '  +2.17%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').Value = '91.44'
$ws.Range('E45').Value = '  +0.90%  '
$ws.Range('E46').Value = '  -10.67%  '
$ws.Range('E47').Value = '  -7.00%  '
$ws.Range('D48').Value = '8.34'
$ws.Range('E48').Value = '  +4.56%  '
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('D50').Value = '2.595.92'
$ws.Range('E50').Value = '  +2.09%  '
$ws.Range('D51').Value = '14.47'
$ws.Range('E51').Value = '  +3.62%  '

# Restore default styling on the cells we touched above so no stray number-format
# survives on the saved worksheet (matches the original "General" formatting).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
